$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Terrain 4" (M column) restriction values for the relevant rows.
$ws.Range("M2").Value = 6
$ws.Range("M5").Value = "6,12"
$ws.Range("M6").Value = "6,12"
$ws.Range("M7").Value = "6,12"
$ws.Range("M10").Value = "6,12"

# Move the active selection (cosmetic, matches author's final cursor position).
$ws.Range("O15").Select()
